$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing A19 timestamp value (tiny correction)
$ws.Range("A19").Value = 45877.79190259259

# Append the new row 20 data
$ws.Range("A20").Value = 45877.83353757377
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 32
$ws.Range("D20").Value = 14.22
$ws.Range("E20").Value = 89.29000000000001
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 6.77
$ws.Range("H20").Value = "ESE"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = "20:00:17"

# Copy the date formatting style from A19 to A20
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
